$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the new test case row (row 4) ---------------------------------
# Copy formatting (border, alignment, number format) from an existing fully
# populated data row (row 2) onto the previously-empty row 4 so the new row
# matches the look of the other test case rows, then set the cell contents.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("A4:G4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A4").Value = 3.2
$ws.Range("B4").Value = "Test deleted to-do item is showing in to-do history page"
$ws.Range("C4").Value = "To test if the to-do item is still showing in the to-do history page after the item is deleted from the to-do page."
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = "The user will be able to see the deleted item in the to-do history page."
$ws.Range("F4").Value = "-"
$ws.Range("G4").Value = "Fail"

# --- Update the view: scroll back to column A and move the selection -------
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("B4").Select() | Out-Null
